# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and mmWave sheets,
# matching the source system's export format (all values stored as text).

$wb = $excel.ActiveWorkbook

function Set-LogRow {
    # NOTE: positional parameters only -- named "-Param value" binding is
    # not reliable against this host, so callers pass args in order.
    param(
        $ws,
        [int]$RowNum,
        [string]$Date,
        [string]$Timestamp,
        [string]$Hour,
        [string]$Location,
        [string]$Value,
        [string]$Status
    )

    # Column A holds a plain "YYYY-MM-DD" string in the source log. Excel's
    # COM layer auto-coerces that pattern into a date serial on assignment,
    # so force the cell to Text first, then strip the number format back to
    # the workbook default (Normal) once the literal text is committed --
    # this keeps the stored value a string without leaving a stray style
    # behind on the cell.
    $cellA = $ws.Cells.Item($RowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $Date
    $cellA.Style = "Normal"

    $ws.Cells.Item($RowNum, 2).Value = $Timestamp
    $ws.Cells.Item($RowNum, 3).Value = $Hour
    $ws.Cells.Item($RowNum, 4).Value = $Location

    # Column E sometimes holds a "NN.N%" reading (Humidity sheet). That
    # pattern also gets auto-coerced (into a percentage number) unless the
    # cell is pre-formatted as Text, so apply the same guard here.
    $cellE = $ws.Cells.Item($RowNum, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $Value
    $cellE.Style = "Normal"

    $ws.Cells.Item($RowNum, 6).Value = $Status
}

# ---- PIR sheet: rows 265-277 (dimension A1:F264 -> A1:F277) ----
$pirWs = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(265, "2026-01-30","17:29:03","17:00","Bathroom","No Motion","Inactive"),
    @(266, "2026-01-30","17:29:04","17:00","Bathroom","No Motion","Inactive"),
    @(267, "2026-01-30","17:29:09","17:00","Bathroom","No Motion","Inactive"),
    @(268, "2026-01-30","17:29:14","17:00","Bathroom","No Motion","Inactive"),
    @(269, "2026-01-30","17:29:19","17:00","Bathroom","No Motion","Inactive"),
    @(270, "2026-01-30","17:29:24","17:00","Bathroom","No Motion","Inactive"),
    @(271, "2026-01-30","17:29:29","17:00","Bathroom","No Motion","Inactive"),
    @(272, "2026-01-30","17:29:34","17:00","Bathroom","No Motion","Inactive"),
    @(273, "2026-01-30","17:29:39","17:00","Bathroom","No Motion","Inactive"),
    @(274, "2026-01-30","17:29:44","17:00","Bathroom","No Motion","Inactive"),
    @(275, "2026-01-30","17:29:49","17:00","Bathroom","No Motion","Inactive"),
    @(276, "2026-01-30","17:29:54","17:00","Bathroom","No Motion","Inactive"),
    @(277, "2026-01-30","17:29:59","17:00","Bathroom","No Motion","Inactive")
)
foreach ($r in $pirRows) {
    Set-LogRow $pirWs $r[0] $r[1] $r[2] $r[3] $r[4] $r[5] $r[6]
}

# ---- Humidity sheet: rows 183-193 (dimension A1:F182 -> A1:F193) ----
$humidityWs = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(183, "2026-01-30","17:29:04","17:00","Bathroom","87.3%","Active"),
    @(184, "2026-01-30","17:29:05","17:00","Bathroom","86.3%","Active"),
    @(185, "2026-01-30","17:29:09","17:00","Bathroom","87.2%","Active"),
    @(186, "2026-01-30","17:29:14","17:00","Bathroom","87.2%","Active"),
    @(187, "2026-01-30","17:29:19","17:00","Bathroom","87.2%","Active"),
    @(188, "2026-01-30","17:29:24","17:00","Bathroom","86.2%","Active"),
    @(189, "2026-01-30","17:29:30","17:00","Bathroom","87.2%","Active"),
    @(190, "2026-01-30","17:29:45","17:00","Bathroom","86.1%","Active"),
    @(191, "2026-01-30","17:29:50","17:00","Bathroom","87.1%","Active"),
    @(192, "2026-01-30","17:29:55","17:00","Bathroom","86.2%","Active"),
    @(193, "2026-01-30","17:30:00","17:00","Bathroom","87.1%","Active")
)
foreach ($r in $humidityRows) {
    Set-LogRow $humidityWs $r[0] $r[1] $r[2] $r[3] $r[4] $r[5] $r[6]
}

# ---- mmWave sheet: row 52 (dimension A1:F51 -> A1:F52) ----
# (Called directly rather than via a single-row array -- a one-item
# @(@(...)) collapses/flattens on this host instead of staying nested.)
$mmwaveWs = $wb.Worksheets.Item("mmWave")
Set-LogRow $mmwaveWs 52 "2026-01-30" "17:29:56" "17:00" "Living Room" "PRESENCE_DETECTED" "Active"
